$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 270, shifting the existing rows 270-281 down to 274-285.
$ws.Rows("270:273").Insert()

# Common column values shared by every data row in this block.
$mercadoId   = 2
$mercado     = "Comercializadora del Agro de Limarí"
$region      = "Coquimbo"
$codreg      = 4
$tipo        = "Fruta"
$productoId  = 100102
$producto    = "Cítricos"
$categoriaId = 100102005
$categoria   = "Naranja"
$unidad      = "$/bins (400 kilos)"
$origen      = "Provincia de Limarí"
$kgUnidad    = 400

$newRows = @(
    @{ Row=270; Fecha=44448; Variedad="Lane Late";  Calidad="Primera"; Volumen=20; PrecioMin=110000; PrecioMax=120000; PrecioProm=115000; PrecioKg=288 },
    @{ Row=271; Fecha=44448; Variedad="Lane Late";  Calidad="Segunda"; Volumen=20; PrecioMin=80000;  PrecioMax=90000;  PrecioProm=85000;  PrecioKg=212 },
    @{ Row=272; Fecha=44448; Variedad="Navel Late"; Calidad="Primera"; Volumen=20; PrecioMin=110000; PrecioMax=120000; PrecioProm=115000; PrecioKg=288 },
    @{ Row=273; Fecha=44448; Variedad="Navel Late"; Calidad="Segunda"; Volumen=20; PrecioMin=80000;  PrecioMax=90000;  PrecioProm=85000;  PrecioKg=212 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = $mercadoId
    $ws.Cells.Item($row, 2).Value2  = $mercado
    $ws.Cells.Item($row, 3).Value2  = $region
    $ws.Cells.Item($row, 4).Value2  = $r.Fecha
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value2  = $codreg
    $ws.Cells.Item($row, 6).Value2  = $tipo
    $ws.Cells.Item($row, 7).Value2  = $productoId
    $ws.Cells.Item($row, 8).Value2  = $producto
    $ws.Cells.Item($row, 9).Value2  = $categoriaId
    $ws.Cells.Item($row, 10).Value2 = $categoria
    $ws.Cells.Item($row, 11).Value2 = $r.Variedad
    $ws.Cells.Item($row, 12).Value2 = $r.Calidad
    $ws.Cells.Item($row, 13).Value2 = $r.Volumen
    $ws.Cells.Item($row, 14).Value2 = $r.PrecioMin
    $ws.Cells.Item($row, 15).Value2 = $r.PrecioMax
    $ws.Cells.Item($row, 16).Value2 = $r.PrecioProm
    $ws.Cells.Item($row, 17).Value2 = $unidad
    $ws.Cells.Item($row, 18).Value2 = $origen
    $ws.Cells.Item($row, 19).Value2 = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value2 = $kgUnidad
}
